# Update countries & provincias Spain
# Daily data refresh of the COVID "Pais" sheet:
#  - A handful of rows get updated totals (new cases landed for that day).
#  - Three country pairs swap rank (and therefore swap which row shows
#    which country's figures) because one of the pair grew past the other:
#      Belgica / Panama      (rows 39-40)
#      Tanzania / Guyana     (rows 164-165)
#      Timor Oriental / Santa Lucia (rows 202-203)
#  - The "last updated" timestamp footer moves from 01:09 to 02:26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 5028913
$ws.Range("C4").Value = 55345
$ws.Range("D4").Value = 2575970
$ws.Range("E4").Value = 2290208
$ws.Range("G4").Value = 1134
$ws.Range("H4").Value = 162735

# --- Row 21: Argentina ---
$ws.Range("B21").Value = 228195
$ws.Range("C21").Value = 7513
$ws.Range("E21").Value = 124092
$ws.Range("G21").Value = 145
$ws.Range("H21").Value = 4251

# --- Row 27: Canada ---
$ws.Range("B27").Value = 118561
$ws.Range("C27").Value = 374
$ws.Range("D27").Value = 103106
$ws.Range("E27").Value = 6489

# --- Rows 39-40: Panama overtakes Belgica ---
$ws.Range("A39").Value = "Panama"
$ws.Range("B39").Value = 71418
$ws.Range("C39").Value = 1187
$ws.Range("D39").Value = 45658
$ws.Range("E39").Value = 24186
$ws.Range("G39").Value = 21
$ws.Range("H39").Value = 1574

$ws.Range("A40").Value = "Belgica"
$ws.Range("B40").Value = 71158
$ws.Range("C40").Value = 510
$ws.Range("D40").Value = 17661
$ws.Range("E40").Value = 43638
$ws.Range("G40").Value = 7
$ws.Range("H40").Value = 9859

# --- Row 98: Paraguay ---
$ws.Range("B98").Value = 6375
$ws.Range("C98").Value = 315
$ws.Range("D98").Value = 4974
$ws.Range("E98").Value = 1335
$ws.Range("G98").Value = 5
$ws.Range("H98").Value = 66

# --- Row 114: Montenegro ---
$ws.Range("B114").Value = 3480
$ws.Range("C114").Value = 69
$ws.Range("D114").Value = 2178
$ws.Range("E114").Value = 1242
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = 60

# --- Row 131: Surinam ---
$ws.Range("B131").Value = 2096
$ws.Range("C131").Value = 46
$ws.Range("D131").Value = 1446
$ws.Range("E131").Value = 621
$ws.Range("G131").Value = 2
$ws.Range("H131").Value = 29

# --- Row 148: Niger ---
$ws.Range("B148").Value = 1153
$ws.Range("C148").Value = 1
$ws.Range("E148").Value = 27

# --- Row 156: Santo Tome y Principe ---
$ws.Range("D156").Value = 797
$ws.Range("E156").Value = 66

# --- Row 158: Bahamas ---
$ws.Range("B158").Value = 761
$ws.Range("C158").Value = 10
$ws.Range("E158").Value = 656

# --- Rows 164-165: Guyana overtakes Tanzania ---
$ws.Range("A164").Value = "Guyana"
$ws.Range("B164").Value = 538
$ws.Range("C164").Value = 29
$ws.Range("D164").Value = 189
$ws.Range("E164").Value = 327
$ws.Range("H164").Value = 22

$ws.Range("A165").Value = "Tanzania"
$ws.Range("D165").Value = 183
$ws.Range("E165").Value = 305
$ws.Range("H165").Value = 21

# --- Rows 202-203: Santa Lucia overtakes Timor Oriental (tied totals) ---
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Agosto de 2020 a las 02:26"
